# Update slug values on the "3-data" sheet (tags table) to their new
# "...3" variants and move the active selection to C8.
#
# Old slug -> New slug
#   t21               -> t31
#   t22               -> t32
#   basic_geo_db      -> basic_geo_db3
#   nature_geo        -> nature_geo3
#   nature_resources  -> nature_resources3
#   nature_res        -> nature_res3
#   eco_res           -> eco_res3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-data")

# Order matters: new shared-string entries are appended in first-seen
# order, so write them in the same order they appear in the saved file.
$ws.Range("A3").Value = "t31"
$ws.Range("A5").Value = "t32"
$ws.Range("C3").Value = "basic_geo_db3"
$ws.Range("C4").Value = "nature_geo3"
$ws.Range("C5").Value = "nature_resources3"
$ws.Range("C6").Value = "nature_res3"
$ws.Range("C7").Value = "eco_res3"

# Move the selection from A2 to C8, as in the edited workbook.
$ws.Range("C8").Select()
